$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7: new entry (No=4, Date=06-Jan-2022, Application=RPA RLOGIC,
# Task=new comment, % of completion=100%, Status=Completed)
$ws.Range("A7").Value = 4

# Copy the date/percent number formats from the row above so the new
# cells keep the same styling (border + date/percent format) instead of
# falling back to General formatting.
$ws.Range("B6").Copy()
$ws.Range("B7").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("B7").Value2 = 44567

$ws.Range("C7").Value = "RPA RLOGIC"
$ws.Range("D7").Value = "1. Created  a sample P&L report with new logic for the MLR and BLR without expenses"

$ws.Range("E6").Copy()
$ws.Range("E7").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("E7").Value2 = 1

$ws.Range("F7").Value = "Completed"

# Row 8: continuation row (Task=new comment, % of completion=50%, Status=WIP)
$ws.Range("D8").Value = "2. P&L report of the HYD is work in progress"

$ws.Range("E6").Copy()
$ws.Range("E8").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("E8").Value2 = 0.5

$ws.Range("F8").Value = "WIP"

$excel.CutCopyMode = 0

# Update the active selection to D20 as recorded in the saved view
$ws.Range("D20").Select()
